# ------------------------------------------------------------------
# Checkpoint before follow-up message
# Updates the GenX FX signal workbook:
#   - Active Signals sheet: refresh first two rows, add 3 new rows
#   - Summary Dashboard sheet: refresh aggregate stats
#   - Signal History sheet: refresh/reshuffle the 15 history rows
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ====================================================================
# Sheet 1: "Active Signals"
# ====================================================================
$ws1 = $wb.Worksheets.Item("Active Signals")

# --- Stash the two "Signal" cell style templates (BUY=green, SELL=red)
#     in scratch cells *before* any edits happen, so later copy/paste
#     operations can't accidentally chain off a half-edited cell.
$ws1.Range("C2").Copy() | Out-Null
$ws1.Range("Z1").PasteSpecial(-4122) | Out-Null   # Z1 = BUY (green) template
$ws1.Range("C3").Copy() | Out-Null
$ws1.Range("Z2").PasteSpecial(-4122) | Out-Null   # Z2 = SELL (red) template

# --- Prepare formatting for the 3 brand-new rows (4,5,6) by copying the
#     whole-row format from the existing row 2 before any values change.
$ws1.Range("A2:J2").Copy() | Out-Null
$ws1.Range("A4:J4").PasteSpecial(-4122) | Out-Null
$ws1.Range("A5:J5").PasteSpecial(-4122) | Out-Null
$ws1.Range("A6:J6").PasteSpecial(-4122) | Out-Null

# --- Column H holds percentages stored as literal text (e.g. "84.0%");
#     force text format first so Excel doesn't coerce them to numbers.
$ws1.Range("H2:H6").NumberFormat = "@"

# Row 2: XAUCHF / BUY  (style already BUY/green - no restyle needed)
$ws1.Range("Z1").Copy() | Out-Null
$ws1.Range("C2").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(2,1).Value = "2025-07-28 20:03"
$ws1.Cells.Item(2,2).Value = "XAUCHF"
$ws1.Cells.Item(2,3).Value = "BUY"
$ws1.Cells.Item(2,4).Value = 2334.28355
$ws1.Cells.Item(2,5).Value = 2334.27912
$ws1.Cells.Item(2,6).Value = 2334.29131
$ws1.Cells.Item(2,7).Value = 0.06
$ws1.Cells.Item(2,8).Value = "84.0%"
$ws1.Cells.Item(2,9).Value = 1.75
$ws1.Cells.Item(2,10).Value = "Active"

# Row 3: NZDUSD / BUY (was SELL-styled XAUUSD -> needs BUY/green style)
$ws1.Range("Z1").Copy() | Out-Null
$ws1.Range("C3").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(3,1).Value = "2025-07-28 19:51"
$ws1.Cells.Item(3,2).Value = "NZDUSD"
$ws1.Cells.Item(3,3).Value = "BUY"
$ws1.Cells.Item(3,4).Value = 0.59221
$ws1.Cells.Item(3,5).Value = 0.58799
$ws1.Cells.Item(3,6).Value = 0.59656
$ws1.Cells.Item(3,7).Value = 0.02
$ws1.Cells.Item(3,8).Value = "76.0%"
$ws1.Cells.Item(3,9).Value = 1.03
$ws1.Cells.Item(3,10).Value = "Active"

# Row 4 (new): XAUCHF / SELL -> needs red SELL style on column C
$ws1.Range("Z2").Copy() | Out-Null
$ws1.Range("C4").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(4,1).Value = "2025-07-28 19:43"
$ws1.Cells.Item(4,2).Value = "XAUCHF"
$ws1.Cells.Item(4,3).Value = "SELL"
$ws1.Cells.Item(4,4).Value = 2336.548
$ws1.Cells.Item(4,5).Value = 2336.55109
$ws1.Cells.Item(4,6).Value = 2336.54131
$ws1.Cells.Item(4,7).Value = 0.07000000000000001
$ws1.Cells.Item(4,8).Value = "87.0%"
$ws1.Cells.Item(4,9).Value = 2.16
$ws1.Cells.Item(4,10).Value = "Active"

# Row 5 (new): USDJPY / SELL -> red SELL style on column C
$ws1.Range("Z2").Copy() | Out-Null
$ws1.Range("C5").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(5,1).Value = "2025-07-28 20:08"
$ws1.Cells.Item(5,2).Value = "USDJPY"
$ws1.Cells.Item(5,3).Value = "SELL"
$ws1.Cells.Item(5,4).Value = 149.07482
$ws1.Cells.Item(5,5).Value = 149.36232
$ws1.Cells.Item(5,6).Value = 148.34779
$ws1.Cells.Item(5,7).Value = 0.04
$ws1.Cells.Item(5,8).Value = "85.0%"
$ws1.Cells.Item(5,9).Value = 2.53
$ws1.Cells.Item(5,10).Value = "Active"

# Row 6 (new): USDJPY / BUY -> green BUY style on column C
$ws1.Range("Z1").Copy() | Out-Null
$ws1.Range("C6").PasteSpecial(-4122) | Out-Null
$ws1.Cells.Item(6,1).Value = "2025-07-28 20:24"
$ws1.Cells.Item(6,2).Value = "USDJPY"
$ws1.Cells.Item(6,3).Value = "BUY"
$ws1.Cells.Item(6,4).Value = 149.10511
$ws1.Cells.Item(6,5).Value = 148.847
$ws1.Cells.Item(6,6).Value = 150.01508
$ws1.Cells.Item(6,7).Value = 0.09
$ws1.Cells.Item(6,8).Value = "81.0%"
$ws1.Cells.Item(6,9).Value = 3.53
$ws1.Cells.Item(6,10).Value = "Active"

# --- Remove the scratch template cells so they don't leak into the sheet
$ws1.Range("Z1:Z2").Clear() | Out-Null

# ====================================================================
# Sheet 2: "Summary Dashboard"
# ====================================================================
$ws2 = $wb.Worksheets.Item("Summary Dashboard")
$ws2.Cells.Item(4,2).Value = 5

# B7 and B8 hold numeric-looking strings ("85.2%", "2.07"); force text
# format so Excel keeps them as literal text instead of coercing them.
$ws2.Cells.Item(7,2).NumberFormat = "@"
$ws2.Cells.Item(7,2).Value = "85.2%"
$ws2.Cells.Item(8,2).NumberFormat = "@"
$ws2.Cells.Item(8,2).Value = "2.07"
$ws2.Cells.Item(9,2).Value = "2025-07-28 20:00:31"

# ====================================================================
# Sheet 3: "Signal History"
# ====================================================================
$ws3 = $wb.Worksheets.Item("Signal History")

function Set-HistoryRow($ws, $row, $ts, $sym, $sig, $entry, $sl, $tp, $lots, $conf, $rr, $status) {
    $ws.Cells.Item($row,1).Value = $ts
    $ws.Cells.Item($row,2).Value = $sym
    $ws.Cells.Item($row,3).Value = $sig
    $ws.Cells.Item($row,4).Value = $entry
    $ws.Cells.Item($row,5).Value = $sl
    $ws.Cells.Item($row,6).Value = $tp
    $ws.Cells.Item($row,7).Value = $lots
    $ws.Cells.Item($row,8).Value = $conf
    $ws.Cells.Item($row,9).Value = $rr
    $ws.Cells.Item($row,10).Value = $status
}

Set-HistoryRow $ws3 2  "2025-07-28 20:10" "USDJPY" "BUY"  149.43539 149.15827 150.33336 0.04 0.83 3.24 "Filled"
Set-HistoryRow $ws3 3  "2025-07-28 20:10" "XAUGBP" "SELL" 2105.39076 2105.39337 2105.38149 0.05 0.9  3.55 "Filled"
Set-HistoryRow $ws3 4  "2025-07-28 20:03" "XAUCHF" "BUY"  2334.28355 2334.27912 2334.29131 0.06 0.84 1.75 "Active"
Set-HistoryRow $ws3 5  "2025-07-28 19:51" "NZDUSD" "BUY"  0.59221 0.58799 0.59656 0.02 0.76 1.03 "Active"
Set-HistoryRow $ws3 6  "2025-07-28 20:01" "XAUAUD" "BUY"  4068.58247 4068.57785 4068.5889 0.05 0.91 1.39 "Filled"
Set-HistoryRow $ws3 7  "2025-07-28 19:43" "XAUCHF" "SELL" 2336.548 2336.55109 2336.54131 0.07000000000000001 0.87 2.16 "Active"
Set-HistoryRow $ws3 8  "2025-07-28 19:41" "EURUSD" "BUY"  1.10507 1.1005 1.10936 0.03 0.83 0.9399999999999999 "Pending"
Set-HistoryRow $ws3 9  "2025-07-28 20:19" "XAUUSD" "BUY"  2644.48224 2644.47753 2644.48859 0.05 0.95 1.35 "Filled"
Set-HistoryRow $ws3 10 "2025-07-28 19:44" "NZDUSD" "SELL" 0.58648 0.58863 0.58035 0.02 0.9399999999999999 2.85 "Filled"
Set-HistoryRow $ws3 11 "2025-07-28 20:08" "USDJPY" "SELL" 149.07482 149.36232 148.34779 0.04 0.85 2.53 "Active"
Set-HistoryRow $ws3 12 "2025-07-28 20:25" "USDCAD" "SELL" 1.36369 1.36737 1.35429 0.09 0.77 2.56 "Filled"
Set-HistoryRow $ws3 13 "2025-07-28 20:24" "USDJPY" "BUY"  149.10511 148.847 150.01508 0.09 0.8100000000000001 3.53 "Active"
Set-HistoryRow $ws3 14 "2025-07-28 19:51" "XAUEUR" "SELL" 2422.95788 2422.96252 2422.95307 0.07000000000000001 0.84 1.04 "Pending"
Set-HistoryRow $ws3 15 "2025-07-28 19:57" "XAUUSD" "SELL" 2649.17888 2649.18361 2649.17293 0.06 0.78 1.26 "Filled"
Set-HistoryRow $ws3 16 "2025-07-28 20:21" "USDCHF" "BUY"  0.88436 0.87957 0.89366 0.07000000000000001 0.9  1.94 "Filled"
